# daily auto push: 2026-02-19 07:14 UTC
# Insert a new data row at row 841 (2026/02/19, 木, 14, 201), pushing the
# existing rows 841-882 down to 842-883. This mirrors the upstream diff,
# which shows the whole 2026/12/29.. block shifting down by one row and a
# brand-new 2026/02/19 row appearing at the top of that block, with the
# sheet's used range growing from D882 to D883.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 841..882 down to 842..883, opening up a blank row 841.
$ws.Rows(841).Insert()

# Column A holds dates formatted as plain text (e.g. "2026/12/29"), not
# real date serials. Format as Text before assigning so COM doesn't
# auto-coerce the string into a date value, then clear the format again
# so the cell ends up with the same (default) style as its neighbours.
$ws.Range("A841").NumberFormat = "@"
$ws.Range("A841").Value = "2026/02/19"
$ws.Range("A841").ClearFormats()

$ws.Range("B841").Value = "木"
$ws.Range("C841").Value = 14
$ws.Range("D841").Value = 201
